# Weekly update: insert 3 new rows of Kiwi price data (newest week) at the
# top of the data block (rows 849-851), shifting the existing rows down.
# This mirrors the commit "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows above the current row 849, pushing the old
# rows 849..874 down to 852..877 (dimension grows from T874 to T877).
$ws.Rows(849).Insert()
$ws.Rows(850).Insert()
$ws.Rows(851).Insert()

# New row 849: Especial, $/caja 18 kilos, Provincia de Curicó
$ws.Range("A849").Value = 9
$ws.Range("B849").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C849").Value = "Metropolitana"
$ws.Range("D849").Value = 45075
$ws.Range("E849").Value = 13
$ws.Range("F849").Value = "Fruta"
$ws.Range("G849").Value = 100101
$ws.Range("H849").Value = "Berries"
$ws.Range("I849").Value = 100101007
$ws.Range("J849").Value = "Kiwi"
$ws.Range("K849").Value = "Hayward"
$ws.Range("L849").Value = "Especial"
$ws.Range("M849").Value = 200
$ws.Range("N849").Value = 14400
$ws.Range("O849").Value = 14400
$ws.Range("P849").Value = 14400
$ws.Range("Q849").Value = "$/caja 18 kilos"
$ws.Range("R849").Value = "Provincia de Curicó"
$ws.Range("S849").Value = 800
$ws.Range("T849").Value = 18

# New row 850: Primera, $/caja 18 kilos, Provincia de Curicó
$ws.Range("A850").Value = 9
$ws.Range("B850").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C850").Value = "Metropolitana"
$ws.Range("D850").Value = 45075
$ws.Range("E850").Value = 13
$ws.Range("F850").Value = "Fruta"
$ws.Range("G850").Value = 100101
$ws.Range("H850").Value = "Berries"
$ws.Range("I850").Value = 100101007
$ws.Range("J850").Value = "Kiwi"
$ws.Range("K850").Value = "Hayward"
$ws.Range("L850").Value = "Primera"
$ws.Range("M850").Value = 220
$ws.Range("N850").Value = 10800
$ws.Range("O850").Value = 10800
$ws.Range("P850").Value = 10800
$ws.Range("Q850").Value = "$/caja 18 kilos"
$ws.Range("R850").Value = "Provincia de Curicó"
$ws.Range("S850").Value = 600
$ws.Range("T850").Value = 18

# New row 851: Segunda, $/caja 18 kilos, Provincia de Curicó
$ws.Range("A851").Value = 9
$ws.Range("B851").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C851").Value = "Metropolitana"
$ws.Range("D851").Value = 45075
$ws.Range("E851").Value = 13
$ws.Range("F851").Value = "Fruta"
$ws.Range("G851").Value = 100101
$ws.Range("H851").Value = "Berries"
$ws.Range("I851").Value = 100101007
$ws.Range("J851").Value = "Kiwi"
$ws.Range("K851").Value = "Hayward"
$ws.Range("L851").Value = "Segunda"
$ws.Range("M851").Value = 180
$ws.Range("N851").Value = 7200
$ws.Range("O851").Value = 7200
$ws.Range("P851").Value = 7200
$ws.Range("Q851").Value = "$/caja 18 kilos"
$ws.Range("R851").Value = "Provincia de Curicó"
$ws.Range("S851").Value = 400
$ws.Range("T851").Value = 18
